$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 82

$ws.Range("A81:V81").Copy()
$ws.Range("A82:V82").PasteSpecial(-4122)

$ws.Cells.Item($row, 1).Value = 81
$ws.Cells.Item($row, 2).Value = "bosnia-and-herzegovina"
$ws.Cells.Item($row, 3).Value = "premijer-liga-bih"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45243.75
$ws.Cells.Item($row, 6).Value = "Velez Mostar"
$ws.Cells.Item($row, 7).Value = 3
$ws.Cells.Item($row, 8).Value = "Zrinjski"
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 2.95
$ws.Cells.Item($row, 11).Value = "12/11/2023 07:12"
$ws.Cells.Item($row, 12).Value = 3.15
$ws.Cells.Item($row, 13).Value = "13/11/2023 17:59"
$ws.Cells.Item($row, 14).Value = 3.03
$ws.Cells.Item($row, 15).Value = "12/11/2023 07:12"
$ws.Cells.Item($row, 16).Value = 3.22
$ws.Cells.Item($row, 17).Value = "13/11/2023 17:55"
$ws.Cells.Item($row, 18).Value = 2.24
$ws.Cells.Item($row, 19).Value = "12/11/2023 07:12"
$ws.Cells.Item($row, 20).Value = 2.26
$ws.Cells.Item($row, 21).Value = "13/11/2023 17:59"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/bosnia-and-herzegovina/premijer-liga-bih/velez-mostar-zrinjski/COT7IiwB/"
